# Add a new "Gelatina" article row just above the existing "Pan" row (row 7),
# shifting all subsequent rows down by one, and backfill the previously-empty
# Imagen path for the "click & roll" cigarrillos row (row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artículos")

# Fill in the missing Imagen value for row 6 (78421974 - click & roll)
$ws.Range("N6").Value = "C:\EditaSoft\Imágenes de artículos\78421974.png"

# Insert a new blank row at position 7, pushing everything else down
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the Gelatina product data
$ws.Range("A7").Value = 7790070432537
$ws.Range("B7").Value = "Gelatina"
$ws.Range("C7").Value = "en polvo"
$ws.Range("D7").Value = "sabor frambuesa"
$ws.Range("E7").Value = "Exquisita"
$ws.Range("F7").Value = 40
$ws.Range("G7").Value = "gr."
$ws.Range("H7").Value = "sobre"
$ws.Range("I7").Value = "Gelatinas"
$ws.Range("J7").Value = "Argentina"
$ws.Range("K7").Value = 12
$ws.Range("L7").Value = $false
$ws.Range("M7").Value = $true
$ws.Range("N7").Value = "C:\EditaSoft\Imágenes de artículos\7790070432537.png"
$ws.Range("O7").Value = $true
$ws.Range("P7").Value = $true
